# Insert a new data row (HIV-1 subtype D, Senegal isolate SE365) above the
# current row 5 ("main refs" sheet), shifting the existing rows 5-19 down to
# 6-20, matching the author's "Subtypes D and G brought into extension build"
# update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main refs")

# Insert a blank row before row 5; everything below (old rows 5-19) shifts
# down to rows 6-20.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new reference record.
$ws.Range("A5").Value = "AB485648"
$ws.Range("B5").Value = "HIV-1"
$ws.Range("C5").Value = "M"
$ws.Range("D5").Value = "D"
$ws.Range("E5").Value = "SE365"
$ws.Range("F5").Value = 1990
$ws.Range("G5").Value = "Senegal"
$ws.Range("H5").Value = "Homo sapiens"
$ws.Range("I5").Value = "human"
$ws.Range("J5").Value = ""

# Match the formatting of the rest of the data rows.
$ws.Range("A5").Style = $ws.Range("A6").Style
$ws.Range("B5").Style = $ws.Range("B6").Style
$ws.Range("C5:J5").Style = $ws.Range("C6:J6").Style

$ws.Range("A1:J20").Select()
$ws.Range("J9").Select()
